$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.412.74'
$ws.Range("E2").Value = '  +1.81%  '

$ws.Range("D3").Value = '3.807.83'
$ws.Range("E3").Value = '  +0.87%  '

$ws.Range("D4").Value = '''1.00'

$ws.Range("D5").Value = '''671.20'
$ws.Range("E5").Value = '  +7.46%  '

$ws.Range("D6").Value = '''169.50'
$ws.Range("E6").Value = '  +2.17%  '

$ws.Range("D7").Value = '3.805.98'
$ws.Range("E7").Value = '  +0.82%  '

$ws.Range("D8").Value = '''1.00'
$ws.Range("E8").Value = '  +0.08%  '

$ws.Range("E9").Value = '  +0.71%  '

$ws.Range("E10").Value = '  +0.95%  '

$ws.Range("D11").Value = '''7.08'
$ws.Range("E11").Value = '  +4.95%  '

$ws.Range("D12").Value = '''0.461'
$ws.Range("E12").Value = '  +0.02%  '

$ws.Range("E13").Value = '  -1.24%  '

$ws.Range("D14").Value = '''35.76'
$ws.Range("E14").Value = '  -0.15%  '

$ws.Range("D15").Value = '4.448.35'
$ws.Range("E15").Value = '  +0.88%  '

$ws.Range("D16").Value = '3.808.04'
$ws.Range("E16").Value = '  +0.93%  '

$ws.Range("D17").Value = '70.442.22'
$ws.Range("E17").Value = '  +1.91%  '

$ws.Range("D18").Value = '''17.69'
$ws.Range("E18").Value = '  +0.17%  '

$ws.Range("D19").Value = '''7.23'
$ws.Range("E19").Value = '  +2.20%  '

$ws.Range("E20").Value = '  +0.53%  '

$ws.Range("D21").Value = '''11.47'
$ws.Range("E21").Value = '  +19.76%  '

$ws.Range("D22").Value = '''477.67'
$ws.Range("E22").Value = '  +2.16%  '

$ws.Range("E23").Value = '  +0.87%  '

$ws.Range("D24").Value = '''83.47'
$ws.Range("E24").Value = '  +0.37%  '

$ws.Range("E25").Value = '  -3.70%  '

$ws.Range("D26").Value = '''12.24'
$ws.Range("E26").Value = '  +1.35%  '

$ws.Range("E27").Value = '  +2.31%  '

$ws.Range("D28").Value = '''2.11'
$ws.Range("E28").Value = '  -2.47%  '

$ws.Range("E29").Value = '  +0.06%  '

$ws.Range("D30").Value = '3.959.16'
$ws.Range("E30").Value = '  +0.90%  '

$ws.Range("D31").Value = '''2.85'
$ws.Range("E31").Value = '  +6.79%  '

$ws.Range("E32").Value = '  +2.22%  '

$ws.Range("E33").Value = '  +3.05%  '

$ws.Range("D34").Value = '''29.62'
$ws.Range("E34").Value = '  +2.83%  '

$ws.Range("E35").Value = '  +6.15%  '

$ws.Range("E36").Value = '  +1.55%  '

$ws.Range("E37").Value = '  -0.01%  '

$ws.Range("D38").Value = '3.763.91'
$ws.Range("E38").Value = '  +1.01%  '

$ws.Range("E39").Value = '  +0.72%  '

$ws.Range("E40").Value = '  +0.37%  '

$ws.Range("D41").Value = '''5.97'
$ws.Range("E41").Value = '  +2.47%  '

$ws.Range("D42").Value = '''0.966'
$ws.Range("E42").Value = '  -0.03%  '

$ws.Range("D43").Value = '''1.00'
$ws.Range("E43").Value = '  +0.07%  '

$ws.Range("D44").Value = '''2.12'
$ws.Range("E44").Value = '  +10.93%  '

$ws.Range("E45").Value = '  -0.01%  '

$ws.Range("D46").Value = '''45.76'
$ws.Range("E46").Value = '  +5.77%  '

$ws.Range("D47").Value = '''158.94'
$ws.Range("E47").Value = '  +4.11%  '

$ws.Range("D48").Value = '''48.10'
$ws.Range("E48").Value = '  +2.96%  '

$ws.Range("B49").Value = 'TheGraph'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D49").Value = '''0.301'
$ws.Range("E49").Value = '  +1.11%  '

$ws.Range("B50").Value = 'FLOKI'
$ws.Range("C50").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D50").Value = '''0.000294'
$ws.Range("E50").Value = '  +6.49%  '

$ws.Range("E51").Value = '  +3.63%  '
